$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6686.1333
$ws.Range("J17").Value = 6834.593
$ws.Range("L17").Value = 20503.779
$ws.Range("N17").Value = -20839.779
$ws.Range("H40").Value = 9819
$ws.Range("J40").Value = 4514.2354
$ws.Range("L40").Value = 4514.2354
$ws.Range("N40").Value = -4864.2354
$ws.Range("H69").Value = 7859.1665
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14126
$ws.Range("H72").Value = 7859.1665
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -40632
$ws.Range("H86").Value = 5935.6294
$ws.Range("I86").Value = 5064.9443
$ws.Range("K86").Value = 5064.9443
$ws.Range("M86").Value = -3941.9443
$ws.Range("H88").Value = 3088.394
$ws.Range("I88").Value = 1298
$ws.Range("J88").Value = 3267.4333
$ws.Range("K88").Value = 1298
$ws.Range("L88").Value = 3267.4333
$ws.Range("M88").Value = -892
$ws.Range("N88").Value = -4079.4333
$ws.Range("H89").Value = 5935.6294
$ws.Range("I89").Value = 5064.9443
$ws.Range("K89").Value = 25324.7215
$ws.Range("M89").Value = -19708.7215
$ws.Range("H91").Value = 3088.394
$ws.Range("I91").Value = 1298
$ws.Range("J91").Value = 3267.4333
$ws.Range("K91").Value = 1298
$ws.Range("L91").Value = 3267.4333
$ws.Range("M91").Value = 106
$ws.Range("N91").Value = -6075.433300000001
$ws.Range("H101").Value = 83333830
$ws.Range("J101").Value = 598.3333
$ws.Range("L101").Value = 1794.9999
$ws.Range("N101").Value = -5038.9999
$ws.Range("H106").Value = 50013780
$ws.Range("I106").Value = 125030200
$ws.Range("J106").Value = 2832
$ws.Range("K106").Value = 125030200
$ws.Range("L106").Value = 2832
$ws.Range("M106").Value = -125029569
$ws.Range("N106").Value = -4094
$ws.Range("H116").Value = 10549.765
$ws.Range("I116").Value = 4751.3
$ws.Range("J116").Value = 18833.285
$ws.Range("K116").Value = 4751.3
$ws.Range("L116").Value = 18833.285
$ws.Range("M116").Value = -1309.3
$ws.Range("N116").Value = -25717.285
$ws.Range("H121").Value = 3353.0908
$ws.Range("J121").Value = 3353.0908
$ws.Range("L121").Value = 10059.2724
$ws.Range("N121").Value = -13553.2724
$ws.Range("H131").Value = 4708.027
$ws.Range("I131").Value = 1224.9412
$ws.Range("J131").Value = 7668.65
$ws.Range("K131").Value = 3674.8236
$ws.Range("L131").Value = 23005.95
$ws.Range("M131").Value = 1365.1764
$ws.Range("N131").Value = -33085.95
$ws.Range("H132").Value = 2994.8774
$ws.Range("I132").Value = 3623
$ws.Range("K132").Value = 10869
$ws.Range("M132").Value = -8339
$ws.Range("H137").Value = 64025.035
$ws.Range("I137").Value = 95813.9
$ws.Range("K137").Value = 287441.7
$ws.Range("M137").Value = -284891.7
$ws.Range("H138").Value = 3242.2363
$ws.Range("I138").Value = 2410.7334
$ws.Range("J138").Value = 3554.05
$ws.Range("K138").Value = 7232.2002
$ws.Range("L138").Value = 10662.15
$ws.Range("M138").Value = -2092.2002
$ws.Range("N138").Value = -20942.15
$ws.Range("H139").Value = 111399.78
$ws.Range("J139").Value = 111399.78
$ws.Range("L139").Value = 111399.78
$ws.Range("N139").Value = -121679.78
$ws.Range("H141").Value = 5542
$ws.Range("I141").Value = 5542
$ws.Range("K141").Value = 16626
$ws.Range("M141").Value = -11446
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1238.13
$ws.Range("I32").Value = 1238.13
$ws.Range("K32").Value = 1238.13
$ws.Range("M32").Value = -951.1300000000001
$ws.Range("H45").Value = 4929750
$ws.Range("I45").Value = 8405386
$ws.Range("K45").Value = 8405386
$ws.Range("M45").Value = -8405009
$ws.Range("H122").Value = 4326298.5
$ws.Range("I122").Value = 7521396.5
$ws.Range("J122").Value = 2089729.6
$ws.Range("K122").Value = 22564189.5
$ws.Range("L122").Value = 6269188.800000001
$ws.Range("M122").Value = -22561739.5
$ws.Range("N122").Value = -6274088.800000001
$ws.Range("H132").Value = 21219.115
$ws.Range("I132").Value = 4141.8486
$ws.Range("J132").Value = 50879.633
$ws.Range("K132").Value = 12425.5458
$ws.Range("L132").Value = 152638.899
$ws.Range("M132").Value = -9895.5458
$ws.Range("N132").Value = -157698.899
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H45").Value = 39000
$ws.Range("J45").Value = 39000
$ws.Range("L45").Value = 39000
$ws.Range("N45").Value = -40616
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2025.7778
$ws.Range("I94").Value = 1705
$ws.Range("K94").Value = 1705
$ws.Range("M94").Value = -1254
$ws.Range("H132").Value = 28132.76
$ws.Range("I132").Value = 20501.154
$ws.Range("K132").Value = 61503.462
$ws.Range("M132").Value = -58973.462
$ws.Range("H134").Value = 8906.844
$ws.Range("I134").Value = 7017.263
$ws.Range("J134").Value = 11668.538
$ws.Range("K134").Value = 21051.789
$ws.Range("L134").Value = 35005.614
$ws.Range("M134").Value = -18516.789
$ws.Range("N134").Value = -40075.614
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2807.5
$ws.Range("I3").Value = 2163.5
$ws.Range("J3").Value = 3880.8333
$ws.Range("K3").Value = 6490.5
$ws.Range("L3").Value = 11642.4999
$ws.Range("M3").Value = -6378.5
$ws.Range("N3").Value = -11866.4999
$ws.Range("H5").Value = 36975.785
$ws.Range("I5").Value = 835.2941
$ws.Range("K5").Value = 2505.8823
$ws.Range("M5").Value = -2393.8823
$ws.Range("H40").Value = 39.9375
$ws.Range("I40").Value = 26.625
$ws.Range("J40").Value = 53.25
$ws.Range("K40").Value = 106.5
$ws.Range("L40").Value = 213
$ws.Range("M40").Value = -37.5
$ws.Range("N40").Value = -351
$ws.Range("H55").Value = 62502080
$ws.Range("I55").Value = 142715550
$ws.Range("J55").Value = 113822.22
$ws.Range("K55").Value = 428146650
$ws.Range("L55").Value = 341466.66
$ws.Range("M55").Value = -428146473
$ws.Range("N55").Value = -341820.66
$ws.Range("H80").Value = 1259.6666
$ws.Range("J80").Value = 1263
$ws.Range("L80").Value = 3789
$ws.Range("N80").Value = -5661
$ws.Range("H83").Value = 1259.6666
$ws.Range("J83").Value = 1263
$ws.Range("L83").Value = 11367
$ws.Range("N83").Value = -20727
$ws.Range("H99").Value = 5374.75
$ws.Range("I99").Value = 5374.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 16124.25
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -13878.25
$ws.Range("N99").ClearContents()
$ws.Range("H113").Value = 2679.861
$ws.Range("I113").Value = 4557.6924
$ws.Range("J113").Value = 1618.4783
$ws.Range("K113").Value = 13673.0772
$ws.Range("L113").Value = 4855.4349
$ws.Range("M113").Value = -11503.0772
$ws.Range("N113").Value = -9195.4349
$ws.Range("H133").Value = 3281.75
$ws.Range("I133").Value = 3281.75
$ws.Range("K133").Value = 9845.25
$ws.Range("M133").Value = -4785.25
$ws.Range("H135").Value = 36975.785
$ws.Range("I135").Value = 835.2941
$ws.Range("K135").Value = 7517.6469
$ws.Range("M135").Value = -4982.6469
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10536324
$ws.Range("J70").Value = 14446.454
$ws.Range("L70").Value = 14446.454
$ws.Range("N70").Value = -14986.454
$ws.Range("H73").Value = 10536324
$ws.Range("J73").Value = 14446.454
$ws.Range("L73").Value = 14446.454
$ws.Range("N73").Value = -16318.454
$ws.Range("H80").Value = 1706579
$ws.Range("I80").Value = 2750362.8
$ws.Range("J80").Value = 364571.44
$ws.Range("K80").Value = 2750362.8
$ws.Range("L80").Value = 364571.44
$ws.Range("M80").Value = -2749364.8
$ws.Range("N80").Value = -366567.44
$ws.Range("H83").Value = 1706579
$ws.Range("I83").Value = 2750362.8
$ws.Range("J83").Value = 364571.44
$ws.Range("K83").Value = 13751814
$ws.Range("L83").Value = 1822857.2
$ws.Range("M83").Value = -13746822
$ws.Range("N83").Value = -1832841.2
$ws.Range("H122").Value = 993799.44
$ws.Range("I122").Value = 1486949.1
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 4460847.300000001
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -4458397.300000001
$ws.Range("N122").Value = -27400
$ws.Range("H132").Value = 5488.6323
$ws.Range("I132").Value = 4195.135
$ws.Range("K132").Value = 12585.405
$ws.Range("M132").Value = -10055.405
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 130691.29
$ws.Range("I22").Value = 444943.5
$ws.Range("J22").Value = 4990.4
$ws.Range("K22").Value = 444943.5
$ws.Range("L22").Value = 4990.4
$ws.Range("M22").Value = -444648.5
$ws.Range("N22").Value = -5580.4
$ws.Range("H27").Value = 130691.29
$ws.Range("I27").Value = 444943.5
$ws.Range("J27").Value = 4990.4
$ws.Range("K27").Value = 444943.5
$ws.Range("L27").Value = 4990.4
$ws.Range("M27").Value = -444836.5
$ws.Range("N27").Value = -5204.4
$ws.Range("H46").Value = 5034.227
$ws.Range("I46").Value = 1387.5
$ws.Range("J46").Value = 5398.9
$ws.Range("K46").Value = 1387.5
$ws.Range("L46").Value = 5398.9
$ws.Range("M46").Value = -1199.5
$ws.Range("N46").Value = -5774.9
$ws.Range("H132").Value = 12671.909
$ws.Range("I132").Value = 13540.566
$ws.Range("J132").Value = 3985.3333
$ws.Range("K132").Value = 40621.698
$ws.Range("L132").Value = 11955.9999
$ws.Range("M132").Value = -38091.698
$ws.Range("N132").Value = -17015.9999
$ws.Range("H136").Value = 100409.09
$ws.Range("J136").Value = 6924.875
$ws.Range("L136").Value = 20774.625
$ws.Range("N136").Value = -25874.625
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11876.286
$ws.Range("I62").Value = 37146
$ws.Range("J62").Value = 8843.92
$ws.Range("K62").Value = 37146
$ws.Range("L62").Value = 8843.92
$ws.Range("M62").Value = -36522
$ws.Range("N62").Value = -10091.92
$ws.Range("H65").Value = 11876.286
$ws.Range("I65").Value = 37146
$ws.Range("J65").Value = 8843.92
$ws.Range("K65").Value = 185730
$ws.Range("L65").Value = 44219.6
$ws.Range("M65").Value = -182610
$ws.Range("N65").Value = -50459.6
